$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# Neighborhood row: was "Talk more", now "Categorical"
$ws.Range("B14").Value = "Categorical"

# Rows whose Type column (B) gets filled in with "Categorical"
$categoricalRows = @(23,24,25,26,27,29,30,31,32,33,34,35,37,41,42,43,44,55,57,59,60,62,65,66,67,74,75,76,77,80)
foreach ($r in $categoricalRows) {
    $ws.Range("B$r").Value = "Categorical"
}

# SaleCondition row gets lowercase "categorical"
$ws.Range("B81").Value = "categorical"

# Leave the selection on the last edited cell, matching where editing stopped
$ws.Range("B81").Select()
